# We should be able to run smoke tests in Preprod and Test
#
# This script:
#  - Adds two new "MHRA" style password values to the shared-string table
#    (in the order the target workbook expects: MHRA@12345A first, then
#    MHRA12345A@) by writing them to the PreProdEnv sheet.
#  - Updates the PreProdEnv (sheet2) B2/B5 passwords to "MHRA12345A@" and
#    the B20 password to "MHRA@12345A" (clearing its special formatting).
#  - Widens PreProdEnv column B so the new, longer passwords fit.
#  - Touches PreProdEnv's page setup (portrait) so it gets its own
#    pageSetup entry / relationship id bump, same as the real edit.
#  - Switches the active sheet/tab from TestEnv to PreProdEnv and updates
#    each sheet's remembered selection (TestEnv -> E21, PreProdEnv -> D7).

$wb = $excel.ActiveWorkbook

$wsTest = $wb.Worksheets.Item("TestEnv")
$wsPreProd = $wb.Worksheets.Item("PreProdEnv")

# --- Update PreProdEnv user/password data ---------------------------------
# Set B20 first so "MHRA@12345A" lands in the shared-string table before
# "MHRA12345A@", matching the order new strings were appended in the edit.
$wsPreProd.Range("B20").Style = "Normal"
$wsPreProd.Range("B20").Value = "MHRA@12345A"

$wsPreProd.Range("B2").Value = "MHRA12345A@"
$wsPreProd.Range("B5").Value = "MHRA12345A@"

# Widen column B on PreProdEnv to fit the longer password strings.
$wsPreProd.Columns.Item(2).ColumnWidth = 13.86

# Touch the page setup (orientation stays portrait) so PreProdEnv gets its
# own print setup entry, same as in the target workbook.
$wsPreProd.PageSetup.Orientation = 1

# --- Update the active sheet / selections ----------------------------------
# Select TestEnv first and place its remembered selection on E21 ...
[void]$wsTest.Select()
[void]$wsTest.Range("E21").Select()

# ... then make PreProdEnv the active tab, with its selection on D7.
[void]$wsPreProd.Select()
[void]$wsPreProd.Range("D7").Select()
